$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F29").Value = 74
$ws.Range("G29").Value = 3791.02
$ws.Range("F30").Value = 142
$ws.Range("G30").Value = 4001.56
$ws.Range("B34").Value = 60167.05
$ws.Range("F36").Value = 90
$ws.Range("G36").Value = 17709.3
$ws.Range("F41").Value = 214
$ws.Range("G41").Value = 41278.46
$ws.Range("F52").Value = 54
$ws.Range("G52").Value = 3186
$ws.Range("F53").Value = 34
$ws.Range("G53").Value = 557.9400000000001
$ws.Range("F56").Value = 40
$ws.Range("G56").Value = 892.8
$ws.Range("F58").Value = 77
$ws.Range("G58").Value = 6000.61
$ws.Range("F61").Value = 236
$ws.Range("G61").Value = 61532.28
$ws.Range("B66").Value = 211084.59
$ws.Range("F151").Value = 36
$ws.Range("G151").Value = 4800.96
$ws.Range("B155").Value = 38348.25
$ws.Range("B161").Value = 57756
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 132.88
$ws.Range("F175").Value = 30
$ws.Range("G175").Value = 8700.299999999999
$ws.Range("F178").Value = 105
$ws.Range("G178").Value = 10159.8
$ws.Range("B193").Value = 67659.92
$ws.Range("F212").Value = 70
$ws.Range("G212").Value = 6236.3
$ws.Range("F213").Value = 221
$ws.Range("G213").Value = 27996.28
$ws.Range("F216").Value = 81
$ws.Range("G216").Value = 6018.3
$ws.Range("F217").Value = 49
$ws.Range("G217").Value = 3640.7
$ws.Range("B218").Value = 82617.21000000001
$ws.Range("F222").Value = 958
$ws.Range("G222").Value = 17723
$ws.Range("B229").Value = 29957.07
$ws.Range("F235").Value = 4
$ws.Range("G235").Value = 1148.04
$ws.Range("B240").Value = 14592.69
$ws.Range("F244").Value = 7
$ws.Range("G244").Value = 3902.5
$ws.Range("B248").Value = 5245.23
$ws.Range("F264").Value = 78
$ws.Range("G264").Value = 2717.52
$ws.Range("F278").Value = 39
$ws.Range("G278").Value = 5286.06
$ws.Range("F284").Value = 171
$ws.Range("G284").Value = 8014.77
$ws.Range("F285").Value = 12
$ws.Range("G285").Value = 1331.28
$ws.Range("F287").Value = 58
$ws.Range("G287").Value = 3174.92
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F291").Value = 27
$ws.Range("G291").Value = 2313.36
$ws.Range("B295").Value = 125143.66
$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 78
$ws.Range("G308").Value = 3715.92
$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12
$ws.Range("F321").Value = 4
$ws.Range("G321").Value = 415.64
$ws.Range("F324").Value = 55
$ws.Range("G324").Value = 9423.15
$ws.Range("B328").Value = -2895.47
$ws.Range("F349").Value = 153
$ws.Range("G349").Value = 11418.39
$ws.Range("F351").Value = 216
$ws.Range("G351").Value = 31242.24
$ws.Range("B356").Value = 79439.14
$ws.Range("F361").Value = 244
$ws.Range("G361").Value = 34303.96
$ws.Range("B363").Value = 78813.97
$ws.Range("F366").Value = 64
$ws.Range("G366").Value = 3541.12
$ws.Range("F368").Value = 61
$ws.Range("G368").Value = 1953.22
$ws.Range("F370").Value = 231
$ws.Range("G370").Value = 38343.69
$ws.Range("B372").Value = 64526.62
$ws.Range("F399").Value = 271
$ws.Range("G399").Value = 26625.75
$ws.Range("F403").Value = 73
$ws.Range("G403").Value = 2959.42
$ws.Range("F408").Value = 211
$ws.Range("G408").Value = 3344.35
$ws.Range("F415").Value = 59
$ws.Range("G415").Value = 3227.3
$ws.Range("B417").Value = 174226.98
$ws.Range("F430").Value = 233
$ws.Range("G430").Value = 10783.24
$ws.Range("F437").Value = 2
$ws.Range("G437").Value = 96.72
$ws.Range("B438").Value = 25993.66
$ws.Range("F452").Value = 55
$ws.Range("G452").Value = 14917.1
$ws.Range("B458").Value = 100492.09
$ws.Range("F478").Value = 12
$ws.Range("G478").Value = 2660.88
$ws.Range("B479").Value = 64810
$ws.Range("E479").Value = 291.22
$ws.Range("F479").Value = 0
$ws.Range("G479").Value = 0
$ws.Range("B480").Value = 53319
$ws.Range("E480").Value = 310.64
$ws.Range("F480").Value = -6
$ws.Range("G480").Value = -1643.52
$ws.Range("B482").Value = 2978.64
$ws.Range("B496").Value = 64833
$ws.Range("E496").Value = 34.9
$ws.Range("F496").Value = 88
$ws.Range("G496").Value = 2889.04
$ws.Range("B497").Value = 60025
$ws.Range("E497").Value = 37.22
$ws.Range("F497").Value = -98
$ws.Range("G497").Value = -3217.34
$ws.Range("B506").Value = 60022
$ws.Range("E506").Value = 37.22
$ws.Range("F506").Value = -113
$ws.Range("G506").Value = -3709.79
$ws.Range("B507").Value = 64830
$ws.Range("E507").Value = 34.9
$ws.Range("F507").Value = 85
$ws.Range("G507").Value = 2790.55
$ws.Range("F512").Value = 25
$ws.Range("G512").Value = 2964.5
$ws.Range("B525").Value = 130213.35
$ws.Range("F527").Value = 56
$ws.Range("G527").Value = 1854.16
$ws.Range("F528").Value = 293
$ws.Range("G528").Value = 4646.98
$ws.Range("F530").Value = 23
$ws.Range("G530").Value = 993.14
$ws.Range("F534").Value = 133
$ws.Range("G534").Value = 5820.08
$ws.Range("B535").Value = 25278.51
$ws.Range("F546").Value = 41
$ws.Range("G546").Value = 6219.7
$ws.Range("B556").Value = 51380.58
$ws.Range("F558").Value = 209
$ws.Range("G558").Value = 25466.65
$ws.Range("F560").Value = 30
$ws.Range("G560").Value = 2411.4
$ws.Range("B561").Value = 29956.85
$ws.Range("F565").Value = 18
$ws.Range("G565").Value = 5057.1
$ws.Range("B573").Value = 27955.37
$ws.Range("F605").Value = 188
$ws.Range("G605").Value = 25022.8
$ws.Range("B607").Value = 25427.83
$ws.Range("F609").Value = 20
$ws.Range("G609").Value = 2176.2
$ws.Range("F614").Value = 82
$ws.Range("G614").Value = 11896.56
$ws.Range("F616").Value = 5
$ws.Range("G616").Value = 713.85
$ws.Range("F617").Value = 27
$ws.Range("G617").Value = 1299.24
$ws.Range("F621").Value = 21
$ws.Range("G621").Value = 7934.01
$ws.Range("F622").Value = 490
$ws.Range("G622").Value = 50425.9
$ws.Range("F623").Value = 80
$ws.Range("G623").Value = 41163.2
$ws.Range("F625").Value = 329
$ws.Range("G625").Value = 12117.07
$ws.Range("B628").Value = 214088.66
$ws.Range("F674").Value = 899
$ws.Range("G674").Value = 146635.89
$ws.Range("B680").Value = 147648.44
$ws.Range("B718").Value = 2822955.79
$ws.Range("B719").Value = 2822955.79
